$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append two new rows to the "Logs" sheet ---
$logs.Range("A24").Value = "Wil graag andere maat ontvangen"
$logs.Range("B24").Value = "mailmind.test@zohomail.eu"
$logs.Range("C24").Value = "Hallo, ik heb het product ontvangen maar de maat is niet goed. Kan ik deze ruilen voor een andere maat?`nSent using {0}"
$logs.Range("D24").Value = "Retour / Terugbetaling"
$logs.Range("E24").Value = "Beste klant,`nBedankt voor je bericht. Vervelend om te horen dat de maat niet juist is. Graag helpen we je verder met het ruilen van het product voor een andere maat. Zou je ons kunnen laten weten welke maat je wilt ontvangen en wat je huidige maat is? Zo kunnen we de ruil snel en correct voor je verwerken.`nWe zien je reactie graag tegemoet.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F24").Value = "2025-06-24 21:00:14"
$logs.Range("G24").Value = "Ja"

$logs.Range("A25").Value = "Offerte voor 500 stuks"
$logs.Range("B25").Value = "mailmind.test@zohomail.eu"
$logs.Range("C25").Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$logs.Range("D25").Value = "Offerte / Prijsaanvraag"
$logs.Range("E25").Value = "Beste klant,`nHartelijk dank voor uw interesse in product X. Om u een nauwkeurige offerte te kunnen sturen, hebben we wat aanvullende informatie nodig, zoals eventuele specifieke wensen met betrekking tot het product of de levering. Kunt u ons ook laten weten naar welk adres de producten verzonden moeten worden?`nZodra we deze gegevens hebben ontvangen, zullen we een offerte voor 500 stuks van product X voor u opstellen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F25").Value = "2025-06-24 21:00:19"
$logs.Range("G25").Value = "Ja"

# --- Extend the conditional formatting ranges to cover the new rows,
#     keeping each rule's dxfId/priority/grouping intact ---
$dFcs = $logs.Range("D2:D23").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D25"))
}

$gFcs = $logs.Range("G2:G23").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G25"))
}

# --- Update the Dashboard summary table (re-sorted counts after the new rows) ---
$dash.Range("B2").Value = 8
$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("A5").Value = "IT / Technisch probleem"
$dash.Range("A6").Value = "Bestelling / Levering"
$dash.Range("B6").Value = 3
$dash.Range("A7").Value = "Sollicitatie / Vacature"
$dash.Range("A8").Value = "Productinformatie"
